$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added at the top of the data block (row 99),
# pushing the existing rows 99-154 down to 100-155.
$ws.Rows.Item(99).Insert()

# Populate the newly inserted row 99 with the new record's data.
$ws.Cells.Item(99, 1).Value = 7
$ws.Cells.Item(99, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(99, 3).Value = "Ñuble"
$ws.Cells.Item(99, 4).Value = 44529
$ws.Cells.Item(99, 5).Value = 16
$ws.Cells.Item(99, 6).Value = 100112003
$ws.Cells.Item(99, 7).Value = "Ajo"
$ws.Cells.Item(99, 8).Value = "Chino"
$ws.Cells.Item(99, 9).Value = "Primera"
$ws.Cells.Item(99, 10).Value = 60
$ws.Cells.Item(99, 11).Value = 17000
$ws.Cells.Item(99, 12).Value = 18000
$ws.Cells.Item(99, 13).Value = 17500
$ws.Cells.Item(99, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(99, 15).Value = "China"
$ws.Cells.Item(99, 16).Value = 1750
$ws.Cells.Item(99, 17).Value = 10
$ws.Cells.Item(99, 18).Value = "Hortaliza"
